# "fixed 0 for benefits"
#
# The upstream report gained a new "Diferença" (Valor OP - Valor Benefício)
# value in column E, the old, always-empty "Resultado" column (F) was
# dropped (shifting every column from G..T one slot to the left, into
# F..S), the benefit lookup for row 3 was refined (a real value + a new
# message instead of "nada encontrado"), and the "Falta Injustificada"
# column (now column M after the shift) was reset to 0 for every
# employee.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the obsolete "Resultado" column (F). Excel shifts G:T left into
#    F:S automatically, which is exactly the column remap seen in the diff.
$ws.Columns("F:F").Delete()

# 2) The newly found benefit value for ELIANE TEREZA DA SILVA (row 3).
$ws.Range("D3").Value = 195.8

# 3) Populate the new "Diferença" column (E) = Valor OP (C) - Valor Benefício (D)
#    for every data row, as a plain value (not a live formula).
$lastRow = 19
for ($r = 2; $r -le $lastRow; $r++) {
    $valorOp = $ws.Cells.Item($r, 3).Value2
    $valorBeneficio = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $valorOp - $valorBeneficio
}

# 4) After the column shift, F holds the old "Alerta" lookup message. Row 3
#    now reports a match across multiple dates instead of "nada encontrado".
$ws.Range("F3").Value = "valor localizado em datas diversas"

# 5) Reset "Falta Injustificada" (column M after the shift) to 0 for every
#    employee - the actual benefits fix described in the commit message.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = 0
}
